$wb = $excel.ActiveWorkbook

# --- Sheet 1: FEINmismatch ---
$ws1 = $wb.Worksheets.Item("FEINmismatch")

# Append two new test-run rows (28, 29), copying the formatting of the last
# existing data row (27) so the new cells pick up the same border / wrap /
# number-format style as the rest of the log.
$ws1.Range("A27:E27").Copy()
$ws1.Range("A28:E28").PasteSpecial(-4122)
$ws1.Range("A29:E29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("E28").Value = "Digital Advertising Gross Revenues"
$ws1.Range("E29").Value = "Digital Advertising Gross Revenues"
$ws1.Range("A28").Value = "Pass"
$ws1.Range("C28").Value = "Y"
$ws1.Range("D28").Value = "Existing Liability with Notice/Invoice Number"
$ws1.Range("A29").Value = "Pass"
$ws1.Range("C29").Value = "Y"
$ws1.Range("D29").Value = "New Tax Return Amount Due"

# Refresh the Date column for every logged run (whole test pack re-executed).
$ws1.Range("B2").Value = "Tue Feb 11 19:34:50 EST 2025"
$ws1.Range("B3").Value = "Tue Feb 11 19:35:03 EST 2025"
$ws1.Range("B4").Value = "Tue Feb 11 19:35:15 EST 2025"
$ws1.Range("B5").Value = "Tue Feb 11 19:35:26 EST 2025"
$ws1.Range("B6").Value = "Tue Feb 11 19:35:38 EST 2025"
$ws1.Range("B7").Value = "Tue Feb 11 19:35:50 EST 2025"
$ws1.Range("B8").Value = "Tue Feb 11 19:36:01 EST 2025"
$ws1.Range("B9").Value = "Tue Feb 11 19:36:12 EST 2025"
$ws1.Range("B10").Value = "Tue Feb 11 19:36:23 EST 2025"
$ws1.Range("B11").Value = "Tue Feb 11 19:36:34 EST 2025"
$ws1.Range("B12").Value = "Tue Feb 11 19:36:45 EST 2025"
$ws1.Range("B13").Value = "Tue Feb 11 19:36:56 EST 2025"
$ws1.Range("B14").Value = "Tue Feb 11 19:37:07 EST 2025"
$ws1.Range("B15").Value = "Tue Feb 11 19:37:18 EST 2025"
$ws1.Range("B16").Value = "Tue Feb 11 19:37:29 EST 2025"
$ws1.Range("B17").Value = "Tue Feb 11 19:37:40 EST 2025"
$ws1.Range("B18").Value = "Tue Feb 11 19:37:51 EST 2025"
$ws1.Range("B19").Value = "Tue Feb 11 19:38:02 EST 2025"
$ws1.Range("B20").Value = "Tue Feb 11 19:38:13 EST 2025"
$ws1.Range("B21").Value = "Tue Feb 11 19:38:24 EST 2025"
$ws1.Range("B22").Value = "Tue Feb 11 19:38:35 EST 2025"
$ws1.Range("B23").Value = "Tue Feb 11 19:38:46 EST 2025"
$ws1.Range("B24").Value = "Tue Feb 11 19:38:57 EST 2025"
$ws1.Range("B25").Value = "Tue Feb 11 19:39:08 EST 2025"
$ws1.Range("B26").Value = "Tue Feb 11 19:39:19 EST 2025"
$ws1.Range("B27").Value = "Tue Feb 11 19:39:30 EST 2025"
$ws1.Range("B28").Value = "Tue Feb 11 19:39:41 EST 2025"
$ws1.Range("B29").Value = "Tue Feb 11 19:39:52 EST 2025"

[void]$ws1.Range("E28").Select()

# --- Sheet 2: FEINSSNmismatch ---
$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")

$ws2.Range("A22:E22").Copy()
$ws2.Range("A23:E23").PasteSpecial(-4122)
$ws2.Range("A24:E24").PasteSpecial(-4122)
$ws2.Range("A25:E25").PasteSpecial(-4122)

$ws2.Range("A23").Value = "Pass"
$ws2.Range("C23").Value = "Y"
$ws2.Range("D23").Value = "Existing Liability with Notice/Invoice Number"
$ws2.Range("E23").Value = "PTE Composite"

$ws2.Range("A24").Value = "Pass"
$ws2.Range("C24").Value = "Y"
$ws2.Range("D24").Value = "New Tax Return Amount Due"
$ws2.Range("E24").Value = "IFTA Tax"

$ws2.Range("A25").Value = "Pass"
$ws2.Range("C25").Value = "Y"
$ws2.Range("D25").Value = "New Tax Return Amount Due"
$ws2.Range("E25").Value = "PTE Composite"

$ws2.Range("B2").Value = "Tue Feb 11 19:40:03 EST 2025"
$ws2.Range("B3").Value = "Tue Feb 11 19:40:14 EST 2025"
$ws2.Range("B4").Value = "Tue Feb 11 19:40:25 EST 2025"
$ws2.Range("B5").Value = "Tue Feb 11 19:40:35 EST 2025"
$ws2.Range("B6").Value = "Tue Feb 11 19:40:46 EST 2025"
$ws2.Range("B7").Value = "Tue Feb 11 19:40:56 EST 2025"
$ws2.Range("B8").Value = "Tue Feb 11 19:41:07 EST 2025"
$ws2.Range("B9").Value = "Tue Feb 11 19:41:18 EST 2025"
$ws2.Range("B10").Value = "Tue Feb 11 19:41:28 EST 2025"
$ws2.Range("B11").Value = "Tue Feb 11 19:41:39 EST 2025"
$ws2.Range("B12").Value = "Tue Feb 11 19:41:50 EST 2025"
$ws2.Range("B13").Value = "Tue Feb 11 19:42:00 EST 2025"
$ws2.Range("B14").Value = "Tue Feb 11 19:42:11 EST 2025"
$ws2.Range("B15").Value = "Tue Feb 11 19:42:22 EST 2025"
$ws2.Range("B16").Value = "Tue Feb 11 19:42:32 EST 2025"
$ws2.Range("B17").Value = "Tue Feb 11 19:42:43 EST 2025"
$ws2.Range("B18").Value = "Tue Feb 11 19:42:53 EST 2025"
$ws2.Range("B19").Value = "Tue Feb 11 19:43:04 EST 2025"
$ws2.Range("B20").Value = "Tue Feb 11 19:43:15 EST 2025"
$ws2.Range("B21").Value = "Tue Feb 11 19:43:25 EST 2025"
$ws2.Range("B22").Value = "Tue Feb 11 19:43:36 EST 2025"
$ws2.Range("B23").Value = "Tue Feb 11 19:43:47 EST 2025"
$ws2.Range("B24").Value = "Tue Feb 11 19:43:57 EST 2025"
$ws2.Range("B25").Value = "Tue Feb 11 19:44:08 EST 2025"

$ws2.Range("E24").Select()
